{"js": "// The document contains a title paragraph with a date, followed by a table\n// whose populated rows hold \"AxB=\" style multiplication expressions.\n// Every text value in the document changes (in document order), including\n// two cells that happen to share the same original text (\"480\u00d76=\") but map\n// to two different replacement values. Because of that we must replace by\n// position (document order) rather than by matching old text.\nconst replacements = [\n  { oldText: \"2025-12-09 Tuesday\", newText: \"2025-12-10 Wednesday\" },\n  { oldText: \"735\u00d73=\", newText: \"483\u00d73=\" },\n  { oldText: \"794\u00d72=\", newText: \"213\u00d79=\" },\n  { oldText: \"919\u00d74=\", newText: \"404\u00d79=\" },\n  { oldText: \"545\u00d73=\", newText: \"232\u00d79=\" },\n  { oldText: \"480\u00d76=\", newText: \"292\u00d74=\" },\n  { oldText: \"480\u00d76=\", newText: \"145\u00d77=\" },\n  { oldText: \"749\u00d73=\", newText: \"252\u00d76=\" },\n  { oldText: \"846\u00d76=\", newText: \"147\u00d76=\" },\n  { oldText: \"574\u00d74=\", newText: \"726\u00d78=\" },\n  { oldText: \"482\u00d72=\", newText: \"119\u00d79=\" },\n  { oldText: \"719\u00d74=\", newText: \"533\u00d75=\" },\n  { oldText: \"806\u00d73=\", newText: \"503\u00d77=\" },\n  { oldText: \"776\u00d72=\", newText: \"933\u00d77=\" },\n  { oldText: \"747\u00d76=\", newText: \"757\u00d77=\" },\n  { oldText: \"993\u00d79=\", newText: \"901\u00d73=\" },\n  { oldText: \"354\u00d78=\", newText: \"448\u00d78=\" },\n  { oldText: \"526\u00d79=\", newText: \"521\u00d75=\" },\n  { oldText: \"931\u00d73=\", newText: \"978\u00d77=\" },\n  { oldText: \"337\u00d73=\", newText: \"616\u00d73=\" },\n  { oldText: \"588\u00d79=\", newText: \"249\u00d78=\" },\n  { oldText: \"642\u00d74=\", newText: \"202\u00d78=\" },\n  { oldText: \"101\u00d76=\", newText: \"678\u00d73=\" },\n  { oldText: \"837\u00d73=\", newText: \"118\u00d73=\" },\n  { oldText: \"620\u00d73=\", newText: \"282\u00d74=\" },\n  { oldText: \"515\u00d78=\", newText: \"245\u00d76=\" },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Only paragraphs that actually contain text are part of the mapping\n// (the table also has many empty paragraphs in blank cells).\nconst nonEmpty = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text !== \"\") {\n    nonEmpty.push(para);\n  }\n}\n\nif (nonEmpty.length !== replacements.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" +\n      replacements.length +\n      \" but found \" +\n      nonEmpty.length\n  );\n}\n\nfor (let i = 0; i < nonEmpty.length; i++) {\n  const expectedOld = replacements[i].oldText;\n  const actualOld = nonEmpty[i].text;\n  if (actualOld !== expectedOld) {\n    throw new Error(\n      \"Mismatch at index \" +\n        i +\n        \": expected '\" +\n        expectedOld +\n        \"' but found '\" +\n        actualOld +\n        \"'\"\n    );\n  }\n  nonEmpty[i].insertText(replacements[i].newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The document contains a title paragraph with a date, followed by a table\n# whose populated rows hold \"AxB=\" style multiplication expressions.\n# Every text value in the document changes (in document order), including\n# two cells that happen to share the same original text (\"480x6=\") but map\n# to two different replacement values. Because of that we replace by\n# position (document order) rather than by matching old text.\n$d = $word.ActiveDocument\n\n$oldTexts = @(\n    \"2025-12-09 Tuesday\",\n    \"735\u00d73=\",\n    \"794\u00d72=\",\n    \"919\u00d74=\",\n    \"545\u00d73=\",\n    \"480\u00d76=\",\n    \"480\u00d76=\",\n    \"749\u00d73=\",\n    \"846\u00d76=\",\n    \"574\u00d74=\",\n    \"482\u00d72=\",\n    \"719\u00d74=\",\n    \"806\u00d73=\",\n    \"776\u00d72=\",\n    \"747\u00d76=\",\n    \"993\u00d79=\",\n    \"354\u00d78=\",\n    \"526\u00d79=\",\n    \"931\u00d73=\",\n    \"337\u00d73=\",\n    \"588\u00d79=\",\n    \"642\u00d74=\",\n    \"101\u00d76=\",\n    \"837\u00d73=\",\n    \"620\u00d73=\",\n    \"515\u00d78=\"\n)\n\n$newTexts = @(\n    \"2025-12-10 Wednesday\",\n    \"483\u00d73=\",\n    \"213\u00d79=\",\n    \"404\u00d79=\",\n    \"232\u00d79=\",\n    \"292\u00d74=\",\n    \"145\u00d77=\",\n    \"252\u00d76=\",\n    \"147\u00d76=\",\n    \"726\u00d78=\",\n    \"119\u00d79=\",\n    \"533\u00d75=\",\n    \"503\u00d77=\",\n    \"933\u00d77=\",\n    \"757\u00d77=\",\n    \"901\u00d73=\",\n    \"448\u00d78=\",\n    \"521\u00d75=\",\n    \"978\u00d77=\",\n    \"616\u00d73=\",\n    \"249\u00d78=\",\n    \"202\u00d78=\",\n    \"678\u00d73=\",\n    \"118\u00d73=\",\n    \"282\u00d74=\",\n    \"245\u00d76=\"\n)\n\n# Walk every paragraph in document order and replace the text of the\n# non-empty ones (skipping empty table-cell paragraphs), matching them\n# one-by-one against the expected original values above.\n$idx = 0\n$total = $d.Paragraphs.Count\nfor ($i = 1; $i -le $total; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $range = $para.Range\n    # Exclude trailing paragraph mark / cell-end marks from comparison.\n    $text = $range.Text\n    $trimmed = $text.TrimEnd([char]13, [char]7)\n    if ($trimmed -ne \"\") {\n        if ($trimmed -ne $oldTexts[$idx]) {\n            throw (\"Mismatch at paragraph \" + $i + \": expected '\" + $oldTexts[$idx] + \"' but found '\" + $trimmed + \"'\")\n        }\n        $range.Text = $newTexts[$idx]\n        $idx = $idx + 1\n    }\n}\n\nif ($idx -ne $oldTexts.Length) {\n    throw (\"Expected to replace \" + $oldTexts.Length + \" paragraphs but replaced \" + $idx)\n}\n"}
